# Applies the "Updated cryptos list" data refresh to sheet1.
# D-column numeric-looking values must stay as TEXT (as in the source data),
# so we force NumberFormat to Text ("@") before assigning those values to
# prevent Excel from auto-converting strings like "213.46" into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.935.15"
$ws.Range("E2").Value = "  +1.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.643.04"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.46"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.90"
$ws.Range("E8").Value = "  +3.00%  "

$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("E11").Value = "  -1.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.79"
$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.647.01"
$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.577"
$ws.Range("E14").Value = "  +5.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.08"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.86"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.915.88"
$ws.Range("E17").Value = "  +1.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.71"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  +0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.64"
$ws.Range("E20").Value = "  +1.60%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.14"
$ws.Range("E22").Value = "  +7.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.41"
$ws.Range("E23").Value = "  +1.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.95"
$ws.Range("E25").Value = "  +2.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("E26").Value = "  +1.06%  "

$ws.Range("E27").Value = "  +0.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.73"
$ws.Range("E28").Value = "  +1.29%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0486"
$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  +2.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.425.88"
$ws.Range("E33").Value = "  -2.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("E34").Value = "  +2.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +2.35%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.892"
$ws.Range("E37").Value = "  +2.15%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.927"
$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0168"
$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.559"
$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("E41").Value = "  +2.50%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.26"
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("B44").Value = "mCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.48"
$ws.Range("E44").Value = "  +0.57%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.46"
$ws.Range("E45").Value = "  +3.34%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.82"
$ws.Range("E46").Value = "  +3.37%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.784.57"
$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "89.01"
$ws.Range("E49").Value = "  +1.95%  "

$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.72"
$ws.Range("E51").Value = "  +0.86%  "
